# Updates cryptos list values (price/volume) per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.958.95'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '1.587.23'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''210.23'
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '''0.245'
$ws.Range("E8").Value = '  -0.60%  '
$ws.Range("E9").Value = '  -1.07%  '
$ws.Range("E10").Value = '  -1.44%  '
$ws.Range("E11").Value = '  +2.29%  '
$ws.Range("D12").Value = '1.809.22'
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("D13").Value = '1.593.83'
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("E14").Value = '  -1.40%  '
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("D16").Value = '25.934.61'
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("D20").Value = '''198.35'
$ws.Range("E20").Value = '  +3.74%  '
$ws.Range("E21").Value = '  +0.29%  '
$ws.Range("D22").Value = '''9.17'
$ws.Range("E22").Value = '  -2.11%  '
$ws.Range("E23").Value = '  +0.63%  '
$ws.Range("E24").Value = '  +8.30%  '
$ws.Range("D25").Value = '''142.60'
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").Value = '''0.120'
$ws.Range("E27").Value = '  -8.91%  '
$ws.Range("D28").Value = '''15.00'
$ws.Range("E28").Value = '  -0.58%  '
$ws.Range("D29").Value = '''6.41'
$ws.Range("E29").Value = '  -0.42%  '
$ws.Range("E30").Value = '  +0.28%  '
$ws.Range("E31").Value = '  +0.21%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("E33").Value = '  -3.13%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '''2.37'
$ws.Range("E34").Value = '  +0.75%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").Value = '''1.46'
$ws.Range("E35").Value = '  -2.05%  '
$ws.Range("D36").Value = '1.120.23'
$ws.Range("E36").Value = '  +1.99%  '
$ws.Range("E37").Value = '  +7.38%  '
$ws.Range("E38").Value = '  -0.20%  '
$ws.Range("E39").Value = '  -1.03%  '
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("D41").Value = '''0.486'
$ws.Range("E41").Value = '  -3.87%  '
$ws.Range("D42").Value = '''0.780'
$ws.Range("E42").Value = '  -5.01%  '
$ws.Range("D43").Value = '1.720.64'
$ws.Range("E43").Value = '  +0.10%  '
$ws.Range("D44").Value = '''5.08'
$ws.Range("E44").Value = '  -2.33%  '
$ws.Range("D45").Value = '''91.76'
$ws.Range("E45").Value = '  -2.16%  '
$ws.Range("E46").Value = '  -2.20%  '
$ws.Range("D47").Value = '''53.07'
$ws.Range("E47").Value = '  -0.22%  '
$ws.Range("D48").Value = '''0.0502'
$ws.Range("E48").Value = '  -1.26%  '
$ws.Range("D49").Value = '''0.406'
$ws.Range("E49").Value = '  -0.39%  '
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("D51").Value = '0.0₇0913'
$ws.Range("E51").Value = '  -18.33%  '
